$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data region (A1:B61) before rewriting, keeping header styles intact
$ws.Range("A1:B61").ClearContents()

# Write final sentence/parse-tree data (some sentences removed, others newly parsed)
$ws.Cells.Item(1, 1).Value = 'Sentence'
$ws.Cells.Item(1, 2).Value = 'Parsed'
$ws.Cells.Item(2, 1).Value = 'Sengap paling Dadong Jepun mangiwasin '
$ws.Cells.Item(2, 2).Value = '(K (S (NP (Noun sengap))) (P (NP (NP (Pronoun paling)) (Noun dadong))) (Pel (NP (NP (Noun jepun)) (Noun mangiwasin))))'
$ws.Cells.Item(3, 1).Value = 'Oo Beli mara teka Luh '
$ws.Cells.Item(3, 2).Value = '(K (S (NP (Noun oo))) (P (NP (Noun beli))) (Ket (PP (Prep mara) (NP (NP (Pronoun teka)) (Noun luh)))))'
$ws.Cells.Item(4, 1).Value = 'Men keto suba matine katepuk '
$ws.Cells.Item(4, 2).Value = '(K (S (NP (NP (Noun men)) (Noun keto))) (P (VP (Adv suba) (VP (Verb matine)))) (Pel (NP (Noun katepuk))))'
$ws.Cells.Item(5, 1).Value = 'Uli ditu baana ngintip '
$ws.Cells.Item(5, 2).Value = '(K (S (NP (Noun uli))) (P (VP (Adv ditu) (VP (Verb baana)))) (Pel (NP (Noun ngintip))))'
$ws.Cells.Item(6, 1).Value = 'Tiang mamusuh teken ia '
$ws.Cells.Item(6, 2).Value = '(K (S (NP (Pronoun tiang))) (P (VP (Verb mamusuh))) (Ket (PP (Prep teken) (NP (Pronoun ia)))))'
$ws.Cells.Item(7, 1).Value = 'jani makejang burone sakti '
$ws.Cells.Item(7, 2).Value = '(K (S (NP (Noun jani))) (P (NP (Pronoun makejang))) (Pel (NP (NP (Pronoun burone)) (Adj sakti))))'
$ws.Cells.Item(8, 1).Value = 'Munyinne suba ngarwanang anak len '
$ws.Cells.Item(8, 2).Value = '(K (S (NP (Pronoun munyinne))) (P (VP (Adv suba) (VP (Verb ngarwanang)))) (Pel (NP (NP (Noun anak)) (Pronoun len))))'
$ws.Cells.Item(9, 1).Value = 'Metu rasa sumanangsaya idane '
$ws.Cells.Item(9, 2).Value = '(K (S (NP (Pronoun metu))) (P (NP (Pronoun rasa))) (Ket (PP (Prep sumanangsaya) (NP (Pronoun idane)))))'
$ws.Cells.Item(10, 1).Value = 'Ia tusing ngelah bulu '
$ws.Cells.Item(10, 2).Value = '(K (S (NP (Pronoun ia))) (P (VP (Adv tusing) (VP (Verb ngelah)))) (Pel (NP (Noun bulu))))'
$ws.Cells.Item(11, 1).Value = 'Kaden tiang ada pancabaya'
$ws.Cells.Item(11, 2).Value = '(K (S (NP (Pronoun kaden))) (P (NP (Pronoun tiang))) (Pel (NP (NP (Noun ada)) (Noun pancabaya))))'
$ws.Cells.Item(12, 1).Value = 'Kenginan pajalanipun setata nyumpling '
$ws.Cells.Item(12, 2).Value = '(K (S (NP (Noun kenginan))) (P (NP (Pronoun pajalanipun))) (Pel (VP (Adv setata) (VP (Verb nyumpling)))))'
$ws.Cells.Item(13, 1).Value = 'sajaan idup lantas buin jani '
$ws.Cells.Item(13, 2).Value = '(K (S (NP (Noun sajaan))) (P (VP (Verb idup))) (O (NP (Noun lantas))) (Pel (NP (NP (Noun buin)) (Noun jani))))'
$ws.Cells.Item(14, 1).Value = 'Serati ubuha limang ukud '
$ws.Cells.Item(14, 2).Value = '(K (S (NP (Noun serati))) (P (NP (NP (Noun ubuha)) (Noun limang))) (Pel (NumP (Num ukud))))'
$ws.Cells.Item(15, 1).Value = 'Nah kanggoang embok masih '
$ws.Cells.Item(15, 2).Value = '(K (S (NP (Pronoun nah))) (P (NP (Noun kanggoang))) (Pel (NP (NP (Noun embok)) (Noun masih))))'
$ws.Cells.Item(16, 1).Value = 'Tulungin ja awake jani '
$ws.Cells.Item(16, 2).Value = '(K (S (NP (Noun tulungin))) (P (NP (Noun ja))) (Pel (NP (NP (Pronoun awake)) (Noun jani))))'
$ws.Cells.Item(17, 1).Value = 'Dadong Jepun suud tiwas idupne '
$ws.Cells.Item(17, 2).Value = '(K (S (NP (Noun dadong))) (P (NP (NP (NP (NP (Noun jepun)) (Noun suud)) (Adj tiwas)) (Pronoun idupne))))'
$ws.Cells.Item(18, 1).Value = 'Ia lantas pules ngengkis '
$ws.Cells.Item(18, 2).Value = '(K (S (NP (Pronoun ia))) (P (NP (Noun lantas))) (Pel (NP (NP (Noun pules)) (Noun ngengkis))))'
$ws.Cells.Item(19, 1).Value = 'Bek misi udang gede-gede pagrepe '
$ws.Cells.Item(19, 2).Value = '(K (S (NP (Noun bek))) (P (NP (NP (Noun misi)) (Pronoun udang))) (Pel (NP (NP (Noun gede-gede)) (Noun pagrepe))))'
$ws.Cells.Item(20, 1).Value = 'Mara neked ditu lantas ungkabanga '
$ws.Cells.Item(20, 2).Value = '(K (S (NP (Pronoun mara))) (P (NP (NP (Pronoun neked)) (Pronoun ditu))) (Pel (NP (NP (Noun lantas)) (Noun ungkabanga))))'
$ws.Cells.Item(21, 1).Value = 'Jani batune dadi dasar temuku '
$ws.Cells.Item(21, 2).Value = '(K (S (NP (Noun jani))) (P (NP (NP (Noun batune)) (Pronoun dadi))) (Pel (NP (NP (Noun dasar)) (Noun temuku))))'
$ws.Cells.Item(22, 1).Value = 'Nah mai tugtug jani icang '
$ws.Cells.Item(22, 2).Value = '(K (S (NP (Pronoun nah))) (P (NP (NP (Noun mai)) (Noun tugtug))) (Pel (NP (NP (Noun jani)) (Pronoun icang))))'
$ws.Cells.Item(23, 1).Value = 'Kenken adi Cucu nakonang '
$ws.Cells.Item(23, 2).Value = '(K (S (NP (Pronoun kenken))) (P (NP (Pronoun adi))) (Pel (NP (NP (Noun cucu)) (Noun nakonang))))'
$ws.Cells.Item(24, 1).Value = 'Makejang pada kedek mabriagan '
$ws.Cells.Item(24, 2).Value = '(K (S (NP (Pronoun makejang))) (P (NP (Noun pada))) (Pel (NP (NP (Noun kedek)) (Noun mabriagan))))'
$ws.Cells.Item(25, 1).Value = 'Keto baos Idane '
$ws.Cells.Item(25, 2).Value = '(K (S (NP (Noun keto))) (P (NP (NP (Noun baos)) (Pronoun idane))))'
$ws.Cells.Item(26, 1).Value = 'Koang keto kone buin aduhanne '
$ws.Cells.Item(26, 2).Value = '(K (S (NP (Noun koang))) (P (NP (Noun keto))) (Pel (VP (Adv kone) (VP (Adv buin) (VP (Verb aduhanne))))))'
$ws.Cells.Item(27, 1).Value = 'Papatihe maparab Ki Patih Bandeswarya '
$ws.Cells.Item(27, 2).Value = '(K (S (NP (Noun papatihe))) (P (NP (NP (Noun maparab)) (Noun ki))) (Pel (NP (NP (Noun patih)) (Noun bandeswarya))))'
$ws.Cells.Item(28, 1).Value = 'Nah wake sing ja kengken '
$ws.Cells.Item(28, 2).Value = '(K (S (NP (Pronoun nah))) (P (NP (NP (Pronoun wake)) (Noun sing))) (Pel (NP (NP (Noun ja)) (Pronoun kengken))))'
$ws.Cells.Item(29, 1).Value = 'Inggih titiang makta '
$ws.Cells.Item(29, 2).Value = '(K (S (NP (Noun inggih))) (P (NP (Pronoun titiang))) (Pel (VP (Verb makta))))'
$ws.Cells.Item(30, 1).Value = 'Keto buin timpalne ngorahin '
$ws.Cells.Item(30, 2).Value = '(K (S (NP (Noun keto))) (P (PP (Prep buin) (NP (NP (Pronoun timpalne)) (Noun ngorahin)))))'
$ws.Cells.Item(31, 1).Value = 'Suud nyurat lantas ia mulih '
$ws.Cells.Item(31, 2).Value = '(K (S (NP (Noun suud))) (P (VP (Verb nyurat))) (O (NP (Noun lantas))) (Pel (NP (NP (Pronoun ia)) (Pronoun mulih))))'
$ws.Cells.Item(32, 1).Value = 'Dadong suba Batara Sri '
$ws.Cells.Item(32, 2).Value = '(K (S (NP (Noun dadong))) (P (VP (Adv suba) (VP (Verb batara)))) (Pel (NP (Noun sri))))'
$ws.Cells.Item(33, 1).Value = 'Ada pandita ajaka duang diri '
$ws.Cells.Item(33, 2).Value = '(K (S (NP (Noun ada))) (P (NP (Noun pandita))) (Ket (PP (Prep ajaka) (NP (NP (Noun duang)) (Noun diri)))))'
$ws.Cells.Item(34, 1).Value = 'Siput ane malunan masaut '
$ws.Cells.Item(34, 2).Value = '(K (S (NP (Noun siput))) (P (PP (Prep ane) (NP (NP (Noun malunan)) (Pronoun masaut)))))'
$ws.Cells.Item(35, 1).Value = 'Nyen adan caine Gede '
$ws.Cells.Item(35, 2).Value = '(K (S (NP (Noun nyen))) (P (NP (Noun adan))) (Pel (NP (NP (Pronoun caine)) (Adj gede))))'
$ws.Cells.Item(36, 1).Value = 'Ditu ia mabakti sambila ngacep '
$ws.Cells.Item(36, 2).Value = '(K (S (NP (Pronoun ditu))) (P (NP (NP (Pronoun ia)) (Noun mabakti))) (Pel (NP (NP (Noun sambila)) (Noun ngacep))))'
$ws.Cells.Item(37, 1).Value = 'Dadinne suung ditu kubune ento '
$ws.Cells.Item(37, 2).Value = '(K (S (NP (Noun dadinne))) (P (NP (NP (Noun suung)) (Pronoun ditu))) (Pel (NP (NP (Noun kubune)) (Pronoun ento))))'
$ws.Cells.Item(38, 1).Value = 'Ada anake ngeling mapangenan '
$ws.Cells.Item(38, 2).Value = '(K (S (NP (Noun ada))) (P (NP (Pronoun anake))) (Pel (NP (NP (Pronoun ngeling)) (Noun mapangenan))))'
$ws.Cells.Item(39, 1).Value = 'Tiang naur nika '
$ws.Cells.Item(39, 2).Value = '(K (S (NP (Pronoun tiang))) (P (NP (NP (Noun naur)) (Pronoun nika))))'
$ws.Cells.Item(40, 1).Value = 'Keto munyin memene '
$ws.Cells.Item(40, 2).Value = '(K (S (NP (Noun keto))) (P (NP (NP (Noun munyin)) (Pronoun memene))))'
$ws.Cells.Item(41, 1).Value = 'Keto bikas alune totonan '
$ws.Cells.Item(41, 2).Value = '(K (S (NP (Noun keto))) (P (VP (Verb bikas))) (Pel (NP (NP (Noun alune)) (Pronoun totonan))))'
$ws.Cells.Item(42, 1).Value = 'Masaut bin curike ne muani '
$ws.Cells.Item(42, 2).Value = '(K (S (NP (Pronoun masaut))) (P (NP (NP (Noun bin)) (Noun curike))) (Pel (NP (NP (Pronoun ne)) (Noun muani))))'
$ws.Cells.Item(43, 1).Value = 'Sajaan ento Dong '
$ws.Cells.Item(43, 2).Value = '(K (S (NP (Noun sajaan))) (P (NP (NP (Pronoun ento)) (Pronoun dong))))'
$ws.Cells.Item(44, 1).Value = 'Inggih titiang sairing '
$ws.Cells.Item(44, 2).Value = '(K (S (NP (Noun inggih))) (P (NP (NP (Pronoun titiang)) (Noun sairing))))'
$ws.Cells.Item(45, 1).Value = 'Ratun kedise uli delod pasih '
$ws.Cells.Item(46, 1).Value = 'Mataluh sabilang wai '
$ws.Cells.Item(46, 2).Value = '(K (S (NP (Pronoun mataluh))) (P (NP (NP (Noun sabilang)) (Noun wai))))'
$ws.Cells.Item(45, 2).Value = '(K (S (NP (Noun ratun))) (P (NP (Pronoun kedise))) (Ket (PP (Prep uli) (NP (NP (Noun delod)) (Noun pasih)))))'

# Update selection to match final saved view state
$ws.Range("B55").Select()
